# Applies the 'new form submission inserted at row 2' edit:
#  - existing data rows 2-20 (columns B:S) shift down to rows 3-21
#  - column A (the sequential index) is left untouched positionally
#  - row 2 (columns B:S) receives the new submission content
#  - row 21 is new; A21 is set to 19 (next index)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing company rows down by one row (columns B:S only) ---
# Row 21 <= old row 20 (A示例xxx公司)
$ws.Cells.Item(21, 2).Value = 'A示例xxx公司'
$ws.Cells.Item(21, 3).Value = 'xx区'
$ws.Cells.Item(21, 4).Value = 'xxx事业部'
$ws.Cells.Item(21, 5).Value = 'Java'
$ws.Cells.Item(21, 6).Value = '9:00-18:30'
$ws.Cells.Item(21, 7).Value = '1.5h'
$ws.Cells.Item(21, 8).Value = '135 加班，24 正常；大小周等等'
$ws.Cells.Item(21, 9).Value = '基数 xxxx，比例 xx%'
$ws.Cells.Item(21, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(21, 11).Value = '是否打折，比如 xx%。'
$ws.Cells.Item(21, 12).Value = '工位大小，环境，是否提供设备，设备型号种类。'
$ws.Cells.Item(21, 13).Value = '是否有入职就有，是否有前置条件才有。'
$ws.Cells.Item(21, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(21, 15).Value = ''
$ws.Cells.Item(21, 16).Value = ''
$ws.Cells.Item(21, 17).Value = '2022-01-24 13:11:01'
$ws.Cells.Item(21, 18).Value = ''
$ws.Cells.Item(21, 19).Value = ''

# Row 20 <= old row 19 (满帮)
$ws.Cells.Item(20, 2).Value = '满帮'
$ws.Cells.Item(20, 3).Value = '雨花区万博科技园'
$ws.Cells.Item(20, 4).Value = ''
$ws.Cells.Item(20, 5).Value = 'Java'
$ws.Cells.Item(20, 6).Value = '9:00-18:30'
$ws.Cells.Item(20, 7).Value = '1.5h'
$ws.Cells.Item(20, 8).Value = '看部门，不强制， 周五基本不加，还有每月一天奋斗日（年底算工资）， 据说要取消了'
$ws.Cells.Item(20, 9).Value = '全额8%'
$ws.Cells.Item(20, 10).Value = '上下半年绩效'
$ws.Cells.Item(20, 11).Value = ''
$ws.Cells.Item(20, 12).Value = '联想'
$ws.Cells.Item(20, 13).Value = ''
$ws.Cells.Item(20, 14).Value = ''
$ws.Cells.Item(20, 15).Value = ''
$ws.Cells.Item(20, 16).Value = ''
$ws.Cells.Item(20, 17).Value = '2022-01-24 14:10:47'
$ws.Cells.Item(20, 18).Value = ''
$ws.Cells.Item(20, 19).Value = ''

# Row 19 <= old row 18 (华为)
$ws.Cells.Item(19, 2).Value = '华为'
$ws.Cells.Item(19, 3).Value = '华为南研所'
$ws.Cells.Item(19, 4).Value = ''
$ws.Cells.Item(19, 5).Value = 'Java'
$ws.Cells.Item(19, 6).Value = '9:00'
$ws.Cells.Item(19, 7).Value = '12:00-13:40'
$ws.Cells.Item(19, 8).Value = '看部门情况。好部门：124加班8：30，35正常下班,差部门：天天11点以后'
$ws.Cells.Item(19, 9).Value = '基础工资的5%'
$ws.Cells.Item(19, 10).Value = '看部门盈利情况和个人绩效定'
$ws.Cells.Item(19, 11).Value = '试用期6个月，100%工资不打折'
$ws.Cells.Item(19, 12).Value = '配win台式机+双屏'
$ws.Cells.Item(19, 13).Value = '没签奋斗协议的5天，但一般不给休，第二年可以换成钱。签了的自愿放弃年假了'
$ws.Cells.Item(19, 14).Value = '必须按时打卡'
$ws.Cells.Item(19, 15).Value = ''
$ws.Cells.Item(19, 16).Value = ''
$ws.Cells.Item(19, 17).Value = '2022-01-24 14:17:32'
$ws.Cells.Item(19, 18).Value = ''
$ws.Cells.Item(19, 19).Value = ''

# Row 18 <= old row 17 (新视云)
$ws.Cells.Item(18, 2).Value = '新视云'
$ws.Cells.Item(18, 3).Value = '雨花台'
$ws.Cells.Item(18, 4).Value = ''
$ws.Cells.Item(18, 5).Value = 'Java'
$ws.Cells.Item(18, 6).Value = '9:00-17:30'
$ws.Cells.Item(18, 7).Value = '1h'
$ws.Cells.Item(18, 8).Value = '看部门，业务部门偶尔加班，技术支持部门基本不加班'
$ws.Cells.Item(18, 9).Value = '基数5k，比例8%'
$ws.Cells.Item(18, 10).Value = '固定13薪'
$ws.Cells.Item(18, 11).Value = '3年合同，试用期总共6个月，前三个月8折，后三个月全薪'
$ws.Cells.Item(18, 12).Value = '配笔记本+显示器'
$ws.Cells.Item(18, 13).Value = '5天年假+5天带薪病假（入职自动折算当年年假）'
$ws.Cells.Item(18, 14).Value = '不打卡'
$ws.Cells.Item(18, 15).Value = ''
$ws.Cells.Item(18, 16).Value = ''
$ws.Cells.Item(18, 17).Value = '2022-01-24 14:17:01'
$ws.Cells.Item(18, 18).Value = ''
$ws.Cells.Item(18, 19).Value = ''

# Row 17 <= old row 16 (创维南京分公司)
$ws.Cells.Item(17, 2).Value = '创维南京分公司'
$ws.Cells.Item(17, 3).Value = '雨花云密城'
$ws.Cells.Item(17, 4).Value = 'web后台'
$ws.Cells.Item(17, 5).Value = 'Java'
$ws.Cells.Item(17, 6).Value = '09:30'
$ws.Cells.Item(17, 7).Value = '1.5h'
$ws.Cells.Item(17, 8).NumberFormat = '@'
$ws.Cells.Item(17, 8).Value = '995'
$ws.Cells.Item(17, 9).Value = '工资八折的10%'
$ws.Cells.Item(17, 10).Value = '1个月工资'
$ws.Cells.Item(17, 11).Value = '不打折'
$ws.Cells.Item(17, 12).Value = 'Windows电脑+dell显示器'
$ws.Cells.Item(17, 13).Value = '法定年假'
$ws.Cells.Item(17, 14).Value = '弹性打卡'
$ws.Cells.Item(17, 15).Value = ''
$ws.Cells.Item(17, 16).Value = ''
$ws.Cells.Item(17, 17).Value = '2022-01-24 14:19:34'
$ws.Cells.Item(17, 18).Value = ''
$ws.Cells.Item(17, 19).Value = ''

# Row 16 <= old row 15 (百家云)
$ws.Cells.Item(16, 2).Value = '百家云'
$ws.Cells.Item(16, 3).Value = '雨花台软件谷科创城'
$ws.Cells.Item(16, 4).Value = ''
$ws.Cells.Item(16, 5).Value = 'Java'
$ws.Cells.Item(16, 6).Value = '9:00-18:30'
$ws.Cells.Item(16, 7).Value = '1.5h'
$ws.Cells.Item(16, 8).Value = '周1,2,4正常加班，不想加班也行'
$ws.Cells.Item(16, 9).Value = ''
$ws.Cells.Item(16, 10).Value = ''
$ws.Cells.Item(16, 11).Value = '6个月不打折。'
$ws.Cells.Item(16, 12).Value = 'mac笔记本+小米曲面屏显示器'
$ws.Cells.Item(16, 13).Value = '年假次年一月发放，每满一年+1天'
$ws.Cells.Item(16, 14).Value = '每个月有4次迟到补卡机会，早上9.15之前打卡不算迟到'
$ws.Cells.Item(16, 15).Value = ''
$ws.Cells.Item(16, 16).Value = ''
$ws.Cells.Item(16, 17).Value = '2022-01-24 14:21:22'
$ws.Cells.Item(16, 18).Value = ''
$ws.Cells.Item(16, 19).Value = ''

# Row 15 <= old row 14 (硅基智能)
$ws.Cells.Item(15, 2).Value = '硅基智能'
$ws.Cells.Item(15, 3).Value = '软件大道'
$ws.Cells.Item(15, 4).Value = '创新产品事业群'
$ws.Cells.Item(15, 5).Value = 'Java'
$ws.Cells.Item(15, 6).Value = '9:00-18:30'
$ws.Cells.Item(15, 7).Value = '1.5h'
$ws.Cells.Item(15, 8).Value = '没事到点走，部门氛围卷'
$ws.Cells.Item(15, 9).Value = '基数5500，比例10%'
$ws.Cells.Item(15, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(15, 11).Value = '不打折'
$ws.Cells.Item(15, 12).Value = '网吧工位'
$ws.Cells.Item(15, 13).Value = '满一年才有正常年假，年假次年一月发放（不满一年打折）'
$ws.Cells.Item(15, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(15, 15).Value = ''
$ws.Cells.Item(15, 16).Value = ''
$ws.Cells.Item(15, 17).Value = '2022-01-24 14:25:34'
$ws.Cells.Item(15, 18).Value = ''
$ws.Cells.Item(15, 19).Value = ''

# Row 14 <= old row 13 (南京力方科技有限公司(力方智充))
$ws.Cells.Item(14, 2).Value = '南京力方科技有限公司(力方智充)'
$ws.Cells.Item(14, 3).Value = '雨花台区软件谷科创城'
$ws.Cells.Item(14, 4).Value = '技术部'
$ws.Cells.Item(14, 5).Value = 'Java'
$ws.Cells.Item(14, 6).Value = '9:00-18:00'
$ws.Cells.Item(14, 7).Value = '1.5h'
$ws.Cells.Item(14, 8).Value = '124固定加班到9点'
$ws.Cells.Item(14, 9).Value = '最低，双边合计512'
$ws.Cells.Item(14, 10).Value = '无'
$ws.Cells.Item(14, 11).Value = '三个月，打八折'
$ws.Cells.Item(14, 12).Value = '网吧工位，自带电脑'
$ws.Cells.Item(14, 13).Value = '法定年假'
$ws.Cells.Item(14, 14).Value = '严格打卡，迟打卡扣30，不打卡半天工资'
$ws.Cells.Item(14, 15).Value = ''
$ws.Cells.Item(14, 16).Value = ''
$ws.Cells.Item(14, 17).Value = '2022-01-24 14:29:37'
$ws.Cells.Item(14, 18).Value = ''
$ws.Cells.Item(14, 19).Value = ''

# Row 13 <= old row 12 (零字节)
$ws.Cells.Item(13, 2).Value = '零字节'
$ws.Cells.Item(13, 3).Value = '建邺'
$ws.Cells.Item(13, 4).Value = ''
$ws.Cells.Item(13, 5).Value = 'Go/Rust/JS/TS/产品/运营'
$ws.Cells.Item(13, 6).Value = '9：30-6：30'
$ws.Cells.Item(13, 7).Value = '1.5h'
$ws.Cells.Item(13, 8).Value = '不加班'
$ws.Cells.Item(13, 9).NumberFormat = '@'
$ws.Cells.Item(13, 9).Value = '8%'
$ws.Cells.Item(13, 10).Value = '13薪，每年调薪一次'
$ws.Cells.Item(13, 11).Value = '应届生八折，有工作经验的不打折'
$ws.Cells.Item(13, 12).Value = 'macbook pro（入职满三年电脑转赠给员工），每人配一个显示器（24-32寸）'
$ws.Cells.Item(13, 13).Value = '入职转正就享受年假'
$ws.Cells.Item(13, 14).Value = '飞书打卡'
$ws.Cells.Item(13, 15).Value = '节日红包、年度旅游（21年三亚一周）'
$ws.Cells.Item(13, 16).Value = ''
$ws.Cells.Item(13, 17).Value = '2022-01-24 14:32:45'
$ws.Cells.Item(13, 18).Value = ''
$ws.Cells.Item(13, 19).Value = ''

# Row 12 <= old row 11 (慧资环球)
$ws.Cells.Item(12, 2).Value = '慧资环球'
$ws.Cells.Item(12, 3).Value = '白下（年中搬到河西）'
$ws.Cells.Item(12, 4).Value = '研发中心'
$ws.Cells.Item(12, 5).Value = '.NET/Python etc.'
$ws.Cells.Item(12, 6).Value = '自己安排，满8小时工时就好'
$ws.Cells.Item(12, 7).Value = '自己安排'
$ws.Cells.Item(12, 8).Value = '不加班'
$ws.Cells.Item(12, 9).Value = '全额8%'
$ws.Cells.Item(12, 10).Value = '13薪，每年调薪一次'
$ws.Cells.Item(12, 11).Value = '不打折'
$ws.Cells.Item(12, 12).Value = '一个高配台式机或者一个高配Dell工作站笔记本，两个40寸4K显示器 Processor Intel(R) Core(TM) i9-10980XE CPU @ 3.00GHz 3.00 GHz  128GB RAM (新的台式机配置标准)'
$ws.Cells.Item(12, 13).Value = '10 ~ 20天'
$ws.Cells.Item(12, 14).Value = '完全不打卡'
$ws.Cells.Item(12, 15).Value = '内推VX：Just1n'
$ws.Cells.Item(12, 16).Value = ''
$ws.Cells.Item(12, 17).Value = '2022-01-24 14:35:55'
$ws.Cells.Item(12, 18).Value = ''
$ws.Cells.Item(12, 19).Value = ''

# Row 11 <= old row 10 (南京希音电子商务有限公司)
$ws.Cells.Item(11, 2).Value = '南京希音电子商务有限公司'
$ws.Cells.Item(11, 3).Value = '天溯产业园'
$ws.Cells.Item(11, 4).Value = ''
$ws.Cells.Item(11, 5).Value = '前端'
$ws.Cells.Item(11, 6).Value = '10:00-18:00(到20:00有50补贴)'
$ws.Cells.Item(11, 7).Value = '12:00-13:30'
$ws.Cells.Item(11, 8).Value = '看部门，不强制，有工时排名。'
$ws.Cells.Item(11, 9).Value = '基础工资的8%'
$ws.Cells.Item(11, 10).Value = '看部门盈利情况和个人绩效定'
$ws.Cells.Item(11, 11).Value = '试用期6个月，100%工资不打折'
$ws.Cells.Item(11, 12).Value = '配mac m1+显示器，网吧工作环境，工位挤。'
$ws.Cells.Item(11, 13).Value = '法定年假，可用加班时长来调休'
$ws.Cells.Item(11, 14).Value = '1月3次补卡'
$ws.Cells.Item(11, 15).Value = '抠，舍得给校招生，不舍得给社招生。多余的调休时长换钱200/d'
$ws.Cells.Item(11, 16).Value = ''
$ws.Cells.Item(11, 17).Value = '2022-01-25 01:58:09'
$ws.Cells.Item(11, 18).Value = ''
$ws.Cells.Item(11, 19).Value = ''

# Row 10 <= old row 9 (南京耀多信息技术有限公司)
$ws.Cells.Item(10, 2).Value = '南京耀多信息技术有限公司'
$ws.Cells.Item(10, 3).Value = '江苏南京'
$ws.Cells.Item(10, 4).Value = '技术部'
$ws.Cells.Item(10, 5).Value = 'Android'
$ws.Cells.Item(10, 6).Value = '9:00-18:00'
$ws.Cells.Item(10, 7).Value = '1h'
$ws.Cells.Item(10, 8).Value = '一开始996，后来发不起加班费不给加班了，欠的加班费也不发'
$ws.Cells.Item(10, 9).Value = '最低额度'
$ws.Cells.Item(10, 10).Value = '无'
$ws.Cells.Item(10, 11).Value = '八折'
$ws.Cells.Item(10, 12).Value = '提供笔记本'
$ws.Cells.Item(10, 13).Value = '有'
$ws.Cells.Item(10, 14).Value = '钉钉位置打卡'
$ws.Cells.Item(10, 15).Value = '老板阴晴不定，随意开除员工'
$ws.Cells.Item(10, 16).Value = ''
$ws.Cells.Item(10, 17).Value = '2022-01-25 02:22:42'
$ws.Cells.Item(10, 18).Value = ''
$ws.Cells.Item(10, 19).Value = ''

# Row 9 <= old row 8 (南京叶子科技有限公司)
$ws.Cells.Item(9, 2).Value = '南京叶子科技有限公司'
$ws.Cells.Item(9, 3).Value = 'xx区'
$ws.Cells.Item(9, 4).Value = 'xxx事业部'
$ws.Cells.Item(9, 5).Value = 'Java'
$ws.Cells.Item(9, 6).Value = '9:00-18:30'
$ws.Cells.Item(9, 7).Value = '1.5h'
$ws.Cells.Item(9, 8).Value = '135 加班，24 正常；大小周等等'
$ws.Cells.Item(9, 9).Value = '基数 xxxx，比例 xx%'
$ws.Cells.Item(9, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(9, 11).Value = '是否打折，比如 xx%。'
$ws.Cells.Item(9, 12).Value = '工位大小，环境，是否提供设备，设备型号种类。'
$ws.Cells.Item(9, 13).Value = '是否有入职就有，是否有前置条件才有。'
$ws.Cells.Item(9, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(9, 15).Value = ''
$ws.Cells.Item(9, 16).Value = ''
$ws.Cells.Item(9, 17).Value = '2022-01-25 07:36:35'
$ws.Cells.Item(9, 18).Value = ''
$ws.Cells.Item(9, 19).Value = ''

# Row 8 <= old row 7 (南京伯索网络科技有限公司（PLASO）)
$ws.Cells.Item(8, 2).Value = '南京伯索网络科技有限公司（PLASO）'
$ws.Cells.Item(8, 3).Value = '秦淮区'
$ws.Cells.Item(8, 4).Value = ''
$ws.Cells.Item(8, 5).Value = ''
$ws.Cells.Item(8, 6).Value = '9:00-18:00'
$ws.Cells.Item(8, 7).Value = '1h'
$ws.Cells.Item(8, 8).Value = '124 加班，35 正常；大小周'
$ws.Cells.Item(8, 9).Value = '基数南京底薪，比例 8%'
$ws.Cells.Item(8, 10).Value = '一般无'
$ws.Cells.Item(8, 11).Value = '3个月8折'
$ws.Cells.Item(8, 12).Value = '网吧工位'
$ws.Cells.Item(8, 13).Value = '入职一年后才有，每年加一天'
$ws.Cells.Item(8, 14).Value = '企业微信打卡，每月三次迟到机会'
$ws.Cells.Item(8, 15).Value = ''
$ws.Cells.Item(8, 16).Value = ''
$ws.Cells.Item(8, 17).Value = '2022-02-06 13:26:16'
$ws.Cells.Item(8, 18).Value = ''
$ws.Cells.Item(8, 19).Value = ''

# Row 7 <= old row 6 (江苏长江汇科技有限公司)
$ws.Cells.Item(7, 2).Value = '江苏长江汇科技有限公司'
$ws.Cells.Item(7, 3).Value = '鼓楼区'
$ws.Cells.Item(7, 4).Value = '技术部'
$ws.Cells.Item(7, 5).Value = 'Java'
$ws.Cells.Item(7, 6).Value = '9:00-17:30'
$ws.Cells.Item(7, 7).Value = '1.5h'
$ws.Cells.Item(7, 8).Value = '按需加班'
$ws.Cells.Item(7, 9).Value = '基数 看个人，比例 5%'
$ws.Cells.Item(7, 10).Value = '基数5000 按照绩效或多或少'
$ws.Cells.Item(7, 11).Value = '8折'
$ws.Cells.Item(7, 12).Value = '无隔板工位，配台式电脑，自带电脑每个月有200补贴，持续24个月。'
$ws.Cells.Item(7, 13).Value = ''
$ws.Cells.Item(7, 14).Value = '钉钉严格打卡'
$ws.Cells.Item(7, 15).Value = ''
$ws.Cells.Item(7, 16).Value = ''
$ws.Cells.Item(7, 17).Value = '2022-02-06 13:30:06'
$ws.Cells.Item(7, 18).Value = ''
$ws.Cells.Item(7, 19).Value = ''

# Row 6 <= old row 5 (云账房)
$ws.Cells.Item(6, 2).Value = '云账房'
$ws.Cells.Item(6, 3).Value = '地址：南站绿地之窗'
$ws.Cells.Item(6, 4).Value = '中台'
$ws.Cells.Item(6, 5).Value = 'Python'
$ws.Cells.Item(6, 6).Value = '9:30-18:30'
$ws.Cells.Item(6, 7).Value = '1.5h'
$ws.Cells.Item(6, 8).Value = '一个月平均加班天数3到4天'
$ws.Cells.Item(6, 9).Value = '总薪资 * 0.8 * 0.6 * 10%'
$ws.Cells.Item(6, 10).Value = '去年只发了半个月'
$ws.Cells.Item(6, 11).Value = '三个月，薪资打八折'
$ws.Cells.Item(6, 12).Value = '网吧工位，一个台式主机，两个24寸1080P显示器'
$ws.Cells.Item(6, 13).Value = '五天年假，按入职日期折算'
$ws.Cells.Item(6, 14).Value = '钉钉打卡，每个月五次补卡机会'
$ws.Cells.Item(6, 15).Value = '根据项目仅仅程度不强制加班。晚上加班白给，周末加班调休'
$ws.Cells.Item(6, 16).Value = ''
$ws.Cells.Item(6, 17).Value = '2022-02-06 13:32:01'
$ws.Cells.Item(6, 18).Value = ''
$ws.Cells.Item(6, 19).Value = ''

# Row 5 <= old row 4 (思杰（Citrix）)
$ws.Cells.Item(5, 2).Value = '思杰（Citrix）'
$ws.Cells.Item(5, 3).Value = '江宁九龙湖'
$ws.Cells.Item(5, 4).Value = ''
$ws.Cells.Item(5, 5).Value = ''
$ws.Cells.Item(5, 6).Value = '09:00-17:00'
$ws.Cells.Item(5, 7).Value = ''
$ws.Cells.Item(5, 8).Value = '不加班'
$ws.Cells.Item(5, 9).Value = '公司缴纳双边12%'
$ws.Cells.Item(5, 10).Value = ''
$ws.Cells.Item(5, 11).Value = ''
$ws.Cells.Item(5, 12).Value = '升降桌+工作站（ 32G+1T ）+MacBook +双显示器+超大工位'
$ws.Cells.Item(5, 13).Value = '年假 15 天，入司满 1 年增加 1 天，上限 20 天'
$ws.Cells.Item(5, 14).Value = ''
$ws.Cells.Item(5, 15).Value = ''
$ws.Cells.Item(5, 16).Value = ''
$ws.Cells.Item(5, 17).Value = '2022-02-07 06:38:11'
$ws.Cells.Item(5, 18).Value = ''
$ws.Cells.Item(5, 19).Value = ''

# Row 4 <= old row 3 (南京叶子科技有限公司)
$ws.Cells.Item(4, 2).Value = '南京叶子科技有限公司'
$ws.Cells.Item(4, 3).Value = '江苏省南京市雨花台区'
$ws.Cells.Item(4, 4).Value = '营销IT事业部'
$ws.Cells.Item(4, 5).Value = 'Java'
$ws.Cells.Item(4, 6).Value = '9:00-18:30'
$ws.Cells.Item(4, 7).Value = '1.5h'
$ws.Cells.Item(4, 8).Value = '正常下班双休，视项目进度自行申请加班'
$ws.Cells.Item(4, 9).Value = ''
$ws.Cells.Item(4, 10).Value = '按KPI决定0~2个月'
$ws.Cells.Item(4, 11).Value = ''
$ws.Cells.Item(4, 12).Value = '工位大小1.5 * 1 长桌，台式机i5-8400 + 8g + 可申请硬盘'
$ws.Cells.Item(4, 13).Value = '入职即可'
$ws.Cells.Item(4, 14).Value = 'OPPO自研IM软件TT打卡'
$ws.Cells.Item(4, 15).Value = ''
$ws.Cells.Item(4, 16).Value = ''
$ws.Cells.Item(4, 17).Value = '2022-02-07 06:37:12'
$ws.Cells.Item(4, 18).Value = '比例 10%左右'
$ws.Cells.Item(4, 19).Value = '试用期 3个月，8折，转正返还'

# Row 3 <= old row 2 (南京三百云信息科技有限公司（车300）)
$ws.Cells.Item(3, 2).Value = '南京三百云信息科技有限公司（车300）'
$ws.Cells.Item(3, 3).Value = '鼓楼区'
$ws.Cells.Item(3, 4).Value = ''
$ws.Cells.Item(3, 5).Value = 'Java'
$ws.Cells.Item(3, 6).Value = '9:00-18:00'
$ws.Cells.Item(3, 7).Value = '1.5h'
$ws.Cells.Item(3, 8).Value = '不强制加班，加班换调休，无加班费，年底清零'
$ws.Cells.Item(3, 9).Value = '基数 工资80%，比例 7%'
$ws.Cells.Item(3, 10).Value = '承诺13薪-14；！！第一年无'
$ws.Cells.Item(3, 11).NumberFormat = '@'
$ws.Cells.Item(3, 11).Value = '100%'
$ws.Cells.Item(3, 12).Value = '提供电脑'
$ws.Cells.Item(3, 13).Value = '每两个月发一天'
$ws.Cells.Item(3, 14).Value = '两次补卡，严格准点打卡。迟到可用调休补(起步0.5h)'
$ws.Cells.Item(3, 15).Value = '抠'
$ws.Cells.Item(3, 16).Value = ''
$ws.Cells.Item(3, 17).Value = '2022-02-08 02:13:11'
$ws.Cells.Item(3, 18).Value = ''
$ws.Cells.Item(3, 19).Value = ''

# --- New row 2: new submission content (columns B:S) ---
$ws.Cells.Item(2, 2).Value = 'xxx有限公司（或缩写）'
$ws.Cells.Item(2, 3).Value = 'xx区'
$ws.Cells.Item(2, 4).Value = 'xxx事业部'
$ws.Cells.Item(2, 5).Value = 'Java'
$ws.Cells.Item(2, 6).Value = '9:00-18:30'
$ws.Cells.Item(2, 7).Value = '1.5h'
$ws.Cells.Item(2, 8).Value = '135 加班，24 正常；大小周等等'
$ws.Cells.Item(2, 9).Value = '基数 xxxx，比例 xx%'
$ws.Cells.Item(2, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(2, 11).Value = '是否打折，比如 xx%。'
$ws.Cells.Item(2, 12).Value = '工位大小，环境，是否提供设备，设备型号种类。'
$ws.Cells.Item(2, 13).Value = '是否有入职就有，是否有前置条件才有。'
$ws.Cells.Item(2, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(2, 15).Value = ''
$ws.Cells.Item(2, 16).Value = ''
$ws.Cells.Item(2, 17).Value = '2022-02-10 02:21:13'
$ws.Cells.Item(2, 18).Value = ''
$ws.Cells.Item(2, 19).Value = ''

# --- New row 21: next sequential index value (match column-A style) ---
$ws.Cells.Item(20, 1).Copy()
$ws.Cells.Item(21, 1).PasteSpecial(-4122)
$ws.Cells.Item(21, 1).Value = 19

